$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix error in Device mapping: "MedicalDevice.Product.ProductType" was
# incorrectly placed next to EHDSDevice.modelNumber (row 10); it belongs
# next to EHDSDevice.type (row 12).
$ws.Range("B10").ClearContents()
$ws.Range("B12").Value = "MedicalDevice.Product.ProductType"

# Update the view/selection to match: scroll back to top and select B12.
$ws.Application.ActiveWindow.ScrollRow = 1
[void]$ws.Range("B12").Select()
